# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# This script rewrites the "Estado de Cuenta" worker/period data table on
# Hoja1 to match the refreshed export: a handful of obsolete period rows
# are removed, the remaining rows are regrouped per worker and re-totalled,
# and the summary header values (VALOR MORA / Cant. Trabajadores /
# Cant. Periodos) are refreshed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ---------------------------------------------------------------------------
# 1) Drop the rows that no longer belong to this statement:
#    - CARLOS ANDRES CABARCAS HERNANDEZ (period 2505) disappears entirely.
#    - LEONARDO CARLOS ALDANA UPARELA loses periods 2501-2504 (keeps 2409-2412).
#    Delete bottom-up so earlier row numbers stay valid while deleting.
# ---------------------------------------------------------------------------
$ws.Range("B37:J37").EntireRow.Delete() | Out-Null
$ws.Range("B18:J21").EntireRow.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2) Rewrite the summary block above the table.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 936583
$ws.Range("C13").Value = 8
$ws.Range("F13").Value = 19

# ---------------------------------------------------------------------------
# 3) Rewrite the 22 data rows (16-37), regrouped by worker, ascending by
#    period within each worker.
# ---------------------------------------------------------------------------
$data = @(
    @("CC", "73185744",   "BLADIMIR PEREZ CARABALLO",        "2102", 13325, 908526),
    @("CC", "73185744",   "BLADIMIR PEREZ CARABALLO",        "2103", 36341, 908526),
    @("CC", "73185744",   "BLADIMIR PEREZ CARABALLO",        "2104", 36341, 908526),
    @("CC", "73185744",   "BLADIMIR PEREZ CARABALLO",        "2105", 36341, 908526),
    @("CC", "1149189090", "EMILIO JOSE CRESPO GOMEZ",        "2204", 40001, 1300000),
    @("CC", "1149189090", "EMILIO JOSE CRESPO GOMEZ",        "2205", 40001, 1300000),
    @("CC", "1149189090", "EMILIO JOSE CRESPO GOMEZ",        "2206", 40001, 1300000),
    @("CC", "1149189090", "EMILIO JOSE CRESPO GOMEZ",        "2207", 40001, 1300000),
    @("CC", "1149189090", "EMILIO JOSE CRESPO GOMEZ",        "2208", 40001, 1300000),
    @("CC", "1149189090", "EMILIO JOSE CRESPO GOMEZ",        "2209", 40001, 1300000),
    @("CC", "1149189090", "EMILIO JOSE CRESPO GOMEZ",        "2210", 40001, 1300000),
    @("CC", "1149189090", "EMILIO JOSE CRESPO GOMEZ",        "2211", 40001, 1300000),
    @("CC", "1149189090", "EMILIO JOSE CRESPO GOMEZ",        "2212", 40001, 1300000),
    @("CC", "1007842115", "YORMAN HERNANDEZ MAZA",           "2306", 46400, 1423500),
    @("CC", "73203527",   "YAN CARLOS RAMOS JULIO",          "2401", 52000, 1300000),
    @("CC", "73153389",   "JOSE MIGUEL NARVAEZ FUENTES",     "2401", 52000, 1300000),
    @("CC", "1002344022", "JANIER DAVID VEGA PEREZ",         "2401", 52000, 1300000),
    @("CC", "1013588837", "LEONARDO CARLOS ALDANA UPARELA",  "2409", 56298, 1456000),
    @("CC", "1013588837", "LEONARDO CARLOS ALDANA UPARELA",  "2410", 58240, 1456000),
    @("CC", "1013588837", "LEONARDO CARLOS ALDANA UPARELA",  "2411", 58240, 1456000),
    @("CC", "1013588837", "LEONARDO CARLOS ALDANA UPARELA",  "2412", 58240, 1456000),
    @("CC", "1140835002", "STEPHANY LOPEZ RODRIGUEZ",        "2412", 20808, 2850000)
)

$row = 16
foreach ($rec in $data) {
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $ws.Cells.Item($row, 6).Value = $rec[4]
    $ws.Cells.Item($row, 7).Value = $rec[5]
    $row = $row + 1
}
